$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slovenia Prva Liga")

# Row 82 and Row 83 had their match-data swapped (id, teams, odds, etc.)
$ws.Range("B82").Value = 6814327
$ws.Range("F82").Value = "NS Mura"
$ws.Range("G82").Value = "NK Domzale"
$ws.Range("I82").Value = 3
$ws.Range("J82").Value = "A"
$ws.Range("K82").Value = 2
$ws.Range("L82").Value = 3.3
$ws.Range("M82").Value = 3.4
$ws.Range("N82").Value = 1.909
$ws.Range("O82").Value = 3.4
$ws.Range("Q82").Value = -0.5
$ws.Range("R82").Value = 1.95
$ws.Range("S82").Value = 1.85
$ws.Range("T82").Value = 2.5
$ws.Range("U82").Value = 1.9
$ws.Range("V82").Value = 1.9
$ws.Range("W82").Value = -1
$ws.Range("X82").Value = -1
$ws.Range("Y82").Value = 2.75
$ws.Range("Z82").Value = -1
$ws.Range("AA82").Value = 0.8500000000000001
$ws.Range("AB82").Value = 0.8999999999999999
$ws.Range("AC82").Value = -1

$ws.Range("B83").Value = 6816473
$ws.Range("F83").Value = "NK Bravo"
$ws.Range("G83").Value = "NK Rogaska"
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = "H"
$ws.Range("K83").Value = 1.8
$ws.Range("L83").Value = 3.5
$ws.Range("M83").Value = 4
$ws.Range("N83").Value = 2.05
$ws.Range("O83").Value = 3
$ws.Range("Q83").Value = -0.25
$ws.Range("R83").Value = 1.75
$ws.Range("S83").Value = 2.05
$ws.Range("T83").Value = 2.25
$ws.Range("U83").Value = 1.95
$ws.Range("V83").Value = 1.85
$ws.Range("W83").Value = 1.05
$ws.Range("X83").Value = -1
$ws.Range("Y83").Value = -1
$ws.Range("Z83").Value = 0.75
$ws.Range("AA83").Value = -1
$ws.Range("AB83").Value = -0.5
$ws.Range("AC83").Value = 0.425

# Rows 181-183 updated to reflect newer odds data, and trailing row 184 removed (matches consolidated)
$ws.Range("B181").Value = 7680774
$ws.Range("E181").Value = 45339.45833333334
$ws.Range("F181").Value = "NK Radomlje"
$ws.Range("G181").Value = "Olimpija Ljubljana"
$ws.Range("K181").Value = 4.8
$ws.Range("L181").Value = 3.8
$ws.Range("M181").Value = 1.615
$ws.Range("N181").Value = 5.25
$ws.Range("O181").Value = 3.75
$ws.Range("P181").Value = 1.571
$ws.Range("Q181").Value = 0.75
$ws.Range("R181").Value = 2.025
$ws.Range("S181").Value = 1.775
$ws.Range("T181").Value = 2.5
$ws.Range("U181").Value = 1.95
$ws.Range("V181").Value = 1.85
$ws.Range("W181").Value = 0
$ws.Range("X181").Value = 0
$ws.Range("Y181").Value = 0
$ws.Range("Z181").Value = 0
$ws.Range("AA181").Value = 0

$ws.Range("B182").Value = 7680773
$ws.Range("E182").Value = 45340.375
$ws.Range("F182").Value = "NK Bravo"
$ws.Range("G182").Value = "NK Celje"
$ws.Range("K182").Value = 4.5
$ws.Range("L182").Value = 3.4
$ws.Range("M182").Value = 1.75
$ws.Range("N182").Value = 6
$ws.Range("O182").Value = 3.5
$ws.Range("P182").Value = 1.615
$ws.Range("Q182").Value = 0.75
$ws.Range("R182").Value = 1.9
$ws.Range("S182").Value = 1.9
$ws.Range("T182").Value = 2.5
$ws.Range("U182").Value = 1.95
$ws.Range("V182").Value = 1.85
$ws.Range("W182").Value = 0
$ws.Range("X182").Value = 0
$ws.Range("Y182").Value = 0
$ws.Range("Z182").Value = 0
$ws.Range("AA182").Value = 0

$ws.Range("B183").Value = 7680776
$ws.Range("E183").Value = 45340.45833333334
$ws.Range("F183").Value = "NS Mura"
$ws.Range("G183").Value = "NK Maribor"
$ws.Range("K183").Value = 3.2
$ws.Range("L183").Value = 3.2
$ws.Range("M183").Value = 2.15
$ws.Range("N183").Value = 3.8
$ws.Range("O183").Value = 3.3
$ws.Range("P183").Value = 1.833
$ws.Range("Q183").Value = 0.5
$ws.Range("R183").Value = 1.975
$ws.Range("S183").Value = 1.825
$ws.Range("T183").Value = 2.75
$ws.Range("U183").Value = 1.975
$ws.Range("V183").Value = 1.825
$ws.Range("W183").Value = 0
$ws.Range("X183").Value = 0
$ws.Range("Y183").Value = 0
$ws.Range("Z183").Value = 0
$ws.Range("AA183").Value = 0

# Remove the now-obsolete last row (184)
$ws.Rows(184).Delete()
